$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Djinn Monk', ['Token Creature — Djinn Monk', 'Flying', '2/2'])"
$ws.Range("A3").Value = "('Dragon', ['Token Creature — Dragon', 'Flying', '4/4'])"
$ws.Range("A4").Value = "('Goblin', ['Token Creature — Goblin', '1/1'])"
$ws.Range("A5").Value = "('Morph', ['Creature', '(You can cover a face-down creature with this reminder card.', 'A card with morph can be turned face up any time for its morph cost.)', '2/2'])"
$ws.Range("A6").Value = "('Narset Transcendent Emblem', ['Emblem — Narset', 'Your opponents can’t cast noncreature spells.'])"
$ws.Range("A7").Value = "('Warrior', ['Token Creature — Warrior', '1/1'])"
$ws.Range("A8").Value = "('Zombie', ['Token Creature — Zombie', '2/2'])"
$ws.Range("A9").Value = "('Zombie Horror', ['Token Creature — Zombie Horror', '*/*'])"

$ws.Range("A10:A29").Clear()

Write-Output $ws.UsedRange.Address()
